$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update source values; dependent formulas (E5 = E4*10, E7 = E6*15) recalc automatically
$ws.Range("E4").Value = 82
$ws.Range("E6").Value = 8

# Move the active selection to E7, matching the saved cursor position
$ws.Range("E7").Select()
